$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-14 Tuesday" "2023-11-15 Wednesday"

Replace-Text "77÷7=11, 0" "65÷7=9, 2"
Replace-Text "10÷7=1, 3" "50÷7=7, 1"
Replace-Text "35÷6=5, 5" "36÷2=18, 0"
Replace-Text "21÷3=7, 0" "19÷8=2, 3"
Replace-Text "96÷7=13, 5" "10÷4=2, 2"

Replace-Text "86÷4=21, 2" "20÷2=10, 0"
Replace-Text "80÷4=20, 0" "35÷4=8, 3"
Replace-Text "93÷7=13, 2" "41÷8=5, 1"
Replace-Text "56÷3=18, 2" "90÷2=45, 0"
Replace-Text "21÷5=4, 1" "32÷7=4, 4"

Replace-Text "28÷8=3, 4" "32÷4=8, 0"
Replace-Text "98÷9=10, 8" "40÷2=20, 0"
Replace-Text "27÷7=3, 6" "14÷8=1, 6"
Replace-Text "20÷8=2, 4" "95÷6=15, 5"
Replace-Text "50÷6=8, 2" "72÷6=12, 0"

Replace-Text "30÷8=3, 6" "48÷3=16, 0"
Replace-Text "71÷4=17, 3" "29÷9=3, 2"
Replace-Text "73÷7=10, 3" "30÷4=7, 2"
Replace-Text "12÷5=2, 2" "95÷4=23, 3"
Replace-Text "37÷5=7, 2" "89÷7=12, 5"

Replace-Text "48÷5=9, 3" "72÷9=8, 0"
Replace-Text "49÷4=12, 1" "95÷4=23, 3"
Replace-Text "50÷4=12, 2" "48÷7=6, 6"
Replace-Text "44÷6=7, 2" "43÷9=4, 7"
Replace-Text "87÷8=10, 7" "59÷8=7, 3"
